$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "PR"
$ws.Range("A9").Value = "HI"
$ws.Range("A10").Select()
